# Append/update timestamp on the "ランサーズ" sheet: rows 2-6, column A
# The original rows were stamped "2026-01-03 12:47:31" and should now read
# "2026-01-03 18:25:48" (commit: Append: 2026-01-03 18:25 JST)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-03 18:25:48"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
